$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.318.70'
$ws.Range("E2").Value = '  -0.12%  '

$ws.Range("D3").Value = '1.842.67'
$ws.Range("E3").Value = '  -1.13%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9991'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.18%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '232.97'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.81%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9987'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.19%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4642'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.78%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2728'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.06%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06271'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.79%  '

$ws.Range("D10").Value = '1.830.06'
$ws.Range("E10").Value = '  -1.82%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07421'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.32%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '16.25'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.34%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.923'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.35%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '83.66'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.48%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6202'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.90%  '

$ws.Range("D16").Value = '30.265.69'
$ws.Range("E16").Value = '  -0.25%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9979'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.27%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '228.07'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.82%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007280'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.16%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.34'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.74%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9993'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.20%  '

$ws.Range("B22").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D22").Value = '2.071.15'
$ws.Range("E22").Value = '  -2.03%  '

$ws.Range("B23").Value = 'Uniswap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.911'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.93%  '

$ws.Range("B24").Value = 'Chainlink'
$ws.Range("C24").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.860'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.06%  '

$ws.Range("B25").Value = 'Cosmos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.175'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.33%  '

$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.38'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.42%  '

$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.79'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.88%  '

$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.868'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.33%  '

$ws.Range("B29").Value = 'Stellar'
$ws.Range("C29").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1035'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.93%  '

$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.369'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.80%  '

$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.076'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.26%  '

$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.811'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.72%  '

$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04841'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.66%  '

$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.142'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.01%  '

$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7106'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.78%  '

$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.693'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.68%  '

$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01877'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.91%  '

$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.655'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.96%  '

$ws.Range("B39").Value = 'TrustWalletToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.8874'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.59%  '

$ws.Range("B40").Value = 'Quant'
$ws.Range("C40").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '104.95'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.84%  '

$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.921'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.89%  '

$ws.Range("B42").Value = 'PaxDollar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.001'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.56%  '

$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.544'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.65%  '

$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4013'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.64%  '

$ws.Range("B45").Value = 'Aptos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.063'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.09%  '

$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '60.52'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.74%  '

$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1195'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.96%  '

$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.655'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.04%  '

$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '33.15'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.75%  '

$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05502'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.21%  '

$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.349'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.55%  '
